$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = -452
$ws.Range("G9").Value = -437
$ws.Range("G11").Value = -437
$ws.Range("G12").Value = -437
